$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (dSF) updates per row, reflecting repulled data / recalculated mean
$values = @{
    3  = -1
    4  = -2
    5  = -1
    6  = -1
    7  = -5
    8  = 1
    9  = -6
    10 = 6
    11 = -4
    12 = 1
    13 = 0
    14 = -6
    15 = -2
    16 = 1
    17 = -1
    18 = 1
    19 = 4
    21 = -8
    22 = -7
    23 = -11
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
